$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
  "3175 The Bays Aged Care Facility Hastings",
  "3563 Embracia Aged Care Reservoir",
  "3600 Belvedere Aged Care Noble Park",
  "3612 BlueCross Glengowrie",
  "3684 Homestyle Aged Care Langford Grange Cranbourne East",
  "3980 Arcare Keysborough Aged Care Keysborough",
  "4075 Ferndale Gardens Aged Care Services Bayswater North",
  "4518 Regis Aged Care Fawkner",
  "AW Window Transport Group Depot North Geelong",
  "Allied Pinnacle Factory Altona North",
  "Bread Solutions Braeside",
  "CS Square Caroline Springs",
  "Cedar Meats Australia Brooklyn",
  "Child's Play Early Learning Centre Tarneit",
  "Comfort Sleep Bedding Co Thomastown",
  "Community Kids Pascoe Vale Early Education Centre Pascoe Vale",
  "Guardian Childcare Caulfield",
  "Hello Fresh Warehouse Ravenhall",
  "Inghams Enterprises Somerville",
  "Lantmannen Unibake Australia Mordialloc",
  "Launch Housing City Edge Crisis Accommodation South Melbourne",
  "Monash Health Kingston Centre South 5",
  "Northern Health Northern Hospital Epping Emergency Department Tier 1B",
  "Northern Health The Northern Hospital Epping",
  "Oceania Meat Processors Laverton North",
  "St Vincents Hospital Emergency Department Melbourne",
  "Target Distribution Centre Truganina",
  "The Robin Hood Inn Drouin West",
  "The Royal Melbourne Hospital Parkville Emergency Department",
  "The Toolshed Bar Private Event Noojee",
  "Turosi Breakwater",
  "Visy Recycling Springvale",
  "Werribee Mercy Hospital Emergency Department",
  "Western Health Footscray Hospital Emergency Department",
  "Western Health Sunshine Hospital Emergency Department"
)

$values = @(16,11,22,35,25,13,17,14,10,12,15,13,10,11,10,22,20,17,24,22,11,10,41,17,16,36,17,46,10,16,10,21,23,10,22)

for ($i = 0; $i -lt $names.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
  $ws.Cells.Item($row, 2).Value = $values[$i]
}
